$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header: volume/issue number and week-covering dates ---
$ws.Range("A8").Value = "Volume 30   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/11/2023  Through  12/17/2023"

# --- Crime Complaints table (rows 14-30): updated weekly/28-day/YTD crime figures ---
$ws.Range("C14").NumberFormat = "#,##0"
$ws.Range("D14").NumberFormat = "#,##0"
$ws.Range("E14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G14").NumberFormat = "#,##0"
$ws.Range("H14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").NumberFormat = "General"
$ws.Range("C16").NumberFormat = "General"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").NumberFormat = "General"
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G27").NumberFormat = "#,##0"
$ws.Range("H27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 7
$ws.Range("J14").Value = 9
$ws.Range("K14").Value = -22.222222222222
$ws.Range("L14").Value = -46.153846153846
$ws.Range("M14").Value = -63.157894736842
$ws.Range("N14").Value = -82.051282051282
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = "0"
$ws.Range("E15").Value = "***.*"
$ws.Range("F15").Value = 4
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 28
$ws.Range("K15").Value = -24.324324324324
$ws.Range("L15").Value = 27.272727272727
$ws.Range("M15").Value = 12
$ws.Range("N15").Value = -63.157894736842
$ws.Range("C16").Value = "0"
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -25
$ws.Range("J16").Value = 126
$ws.Range("K16").Value = 6.349206349206
$ws.Range("L16").Value = -14.102564102564
$ws.Range("M16").Value = -60.588235294117
$ws.Range("N16").Value = -88.025022341376
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -11.111111111111
$ws.Range("F17").Value = 29
$ws.Range("G17").Value = 33
$ws.Range("H17").Value = -12.121212121212
$ws.Range("I17").Value = 343
$ws.Range("J17").Value = 376
$ws.Range("K17").Value = -8.77659574468
$ws.Range("L17").Value = -11.139896373057
$ws.Range("M17").Value = 0.882352941176
$ws.Range("N17").Value = -53.523035230352
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -54.545454545454
$ws.Range("I18").Value = 80
$ws.Range("J18").Value = 123
$ws.Range("K18").Value = -34.959349593495
$ws.Range("L18").Value = -35.483870967741
$ws.Range("M18").Value = -80.19801980198
$ws.Range("N18").Value = -92.727272727272
$ws.Range("C19").Value = 11
$ws.Range("E19").Value = 57.142857142857
$ws.Range("F19").Value = 23
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = -11.538461538461
$ws.Range("I19").Value = 365
$ws.Range("J19").Value = 424
$ws.Range("K19").Value = -13.915094339622
$ws.Range("L19").Value = 4.58452722063
$ws.Range("M19").Value = -43.234836702954
$ws.Range("N19").Value = -90.464994775339
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 35
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = 75
$ws.Range("I20").Value = 223
$ws.Range("J20").Value = 227
$ws.Range("K20").Value = -1.762114537444
$ws.Range("L20").Value = 17.368421052631
$ws.Range("M20").Value = -14.230769230769
$ws.Range("N20").Value = -87.401129943502
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = -6.896551724137
$ws.Range("F21").Value = 106
$ws.Range("G21").Value = 105
$ws.Range("H21").Value = 0.95238095238
$ws.Range("I21").Value = 1180
$ws.Range("J21").Value = 1322
$ws.Range("K21").Value = -10.741301059001
$ws.Range("L21").Value = -4.838709677419
$ws.Range("M21").Value = -41.90054160512
$ws.Range("N21").Value = -86.38985005767
$ws.Range("C24").Value = 38
$ws.Range("D24").Value = 35
$ws.Range("E24").Value = 8.571428571428
$ws.Range("F24").Value = 117
$ws.Range("G24").Value = 108
$ws.Range("H24").Value = 8.333333333333
$ws.Range("I24").Value = 1144
$ws.Range("J24").Value = 1279
$ws.Range("K24").Value = -10.555121188428
$ws.Range("L24").Value = 20.16806722689
$ws.Range("M24").Value = 8.95238095238
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = 14.285714285714
$ws.Range("F25").Value = 63
$ws.Range("G25").Value = 52
$ws.Range("H25").Value = 21.153846153846
$ws.Range("I25").Value = 670
$ws.Range("J25").Value = 524
$ws.Range("K25").Value = 27.862595419847
$ws.Range("L25").Value = 42.553191489361
$ws.Range("M25").Value = -13.659793814433
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = "0"
$ws.Range("E26").Value = "***.*"
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 33.333333333333
$ws.Range("I26").Value = 46
$ws.Range("K26").Value = -11.538461538461
$ws.Range("L26").Value = 12.195121951219
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 150
$ws.Range("I27").Value = 45
$ws.Range("J27").Value = 54
$ws.Range("K27").Value = -16.666666666666
$ws.Range("L27").Value = -8.163265306122
$ws.Range("C28").Value = 2
$ws.Range("F28").Value = 4
$ws.Range("I28").Value = 22
$ws.Range("K28").Value = -29.032258064516
$ws.Range("L28").Value = -59.259259259259
$ws.Range("M28").Value = -62.068965517241
$ws.Range("N28").Value = -85.897435897435
$ws.Range("F29").Value = 3
$ws.Range("I29").Value = 16
$ws.Range("K29").Value = -36
$ws.Range("L29").Value = -62.790697674418
$ws.Range("M29").Value = -66.666666666666
$ws.Range("N29").Value = -88.571428571428

# --- Row 15 label stays "Rape" (shared-string reorder only, no visible text change) ---
